$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(149, 2).Value = 48654
$ws.Cells.Item(149, 5).Value = 38.26
$ws.Cells.Item(149, 6).Value = -1
$ws.Cells.Item(149, 7).Value = -32.02

$ws.Cells.Item(150, 2).Value = 63902
$ws.Cells.Item(150, 5).Value = 34.04
$ws.Cells.Item(150, 6).Value = 2
$ws.Cells.Item(150, 7).Value = 64.04000000000001

$ws.Cells.Item(161, 2).Value = 53925
$ws.Cells.Item(161, 6).Value = 1
$ws.Cells.Item(161, 7).Value = 66.44

$ws.Cells.Item(163, 2).Value = 57756
$ws.Cells.Item(163, 6).Value = -100
$ws.Cells.Item(163, 7).Value = -6644

$ws.Cells.Item(183, 2).Value = 64329
$ws.Cells.Item(183, 5).Value = 128.32
$ws.Cells.Item(183, 6).Value = 6
$ws.Cells.Item(183, 7).Value = 724.14

$ws.Cells.Item(184, 2).Value = 57552
$ws.Cells.Item(184, 5).Value = 136.86
$ws.Cells.Item(184, 6).Value = -5
$ws.Cells.Item(184, 7).Value = -603.45

$ws.Cells.Item(313, 2).Value = 57854
$ws.Cells.Item(313, 6).Value = 2
$ws.Cells.Item(313, 7).Value = 611.6799999999999

$ws.Cells.Item(314, 2).Value = 62997
$ws.Cells.Item(314, 6).Value = 72
$ws.Cells.Item(314, 7).Value = 22020.48

$ws.Cells.Item(316, 2).Value = 57077
$ws.Cells.Item(316, 4).Value = 93.08
$ws.Cells.Item(316, 5).Value = 111.2
$ws.Cells.Item(316, 6).Value = 1
$ws.Cells.Item(316, 7).Value = 93.08

$ws.Cells.Item(317, 2).Value = 61610
$ws.Cells.Item(317, 5).Value = 122.71
$ws.Cells.Item(317, 6).Value = -58
$ws.Cells.Item(317, 7).Value = -5957.18

$ws.Cells.Item(318, 2).Value = 63565
$ws.Cells.Item(318, 4).Value = 102.71
$ws.Cells.Item(318, 5).Value = 109.19
$ws.Cells.Item(318, 6).Value = 60
$ws.Cells.Item(318, 7).Value = 6162.6

$ws.Cells.Item(346, 2).Value = 55373
$ws.Cells.Item(346, 5).Value = 163.62
$ws.Cells.Item(346, 6).Value = -94
$ws.Cells.Item(346, 7).Value = -13562.32

$ws.Cells.Item(347, 2).Value = 63520
$ws.Cells.Item(347, 5).Value = 153.4
$ws.Cells.Item(347, 6).Value = 97
$ws.Cells.Item(347, 7).Value = 13995.16

$ws.Cells.Item(351, 2).Value = 57802
$ws.Cells.Item(351, 5).Value = 162.71
$ws.Cells.Item(351, 6).Value = -79
$ws.Cells.Item(351, 7).Value = -11334.92

$ws.Cells.Item(352, 2).Value = 63531
$ws.Cells.Item(352, 5).Value = 152.53
$ws.Cells.Item(352, 6).Value = 80
$ws.Cells.Item(352, 7).Value = 11478.4

$ws.Cells.Item(355, 2).Value = 63510
$ws.Cells.Item(355, 5).Value = 50.66
$ws.Cells.Item(355, 6).Value = 167
$ws.Cells.Item(355, 7).Value = 7955.88

$ws.Cells.Item(356, 2).Value = 55356
$ws.Cells.Item(356, 5).Value = 54.04
$ws.Cells.Item(356, 6).Value = -158
$ws.Cells.Item(356, 7).Value = -7527.12

$ws.Cells.Item(375, 2).Value = 61605
$ws.Cells.Item(375, 5).Value = 133.78
$ws.Cells.Item(375, 6).Value = -13
$ws.Cells.Item(375, 7).Value = -1455.48

$ws.Cells.Item(376, 2).Value = 63563
$ws.Cells.Item(376, 5).Value = 119.04
$ws.Cells.Item(376, 6).Value = 15
$ws.Cells.Item(376, 7).Value = 1679.4

$ws.Cells.Item(379, 2).Value = 63564
$ws.Cells.Item(379, 5).Value = 137.16
$ws.Cells.Item(379, 6).Value = 57
$ws.Cells.Item(379, 7).Value = 7353.57

$ws.Cells.Item(380, 2).Value = 61608
$ws.Cells.Item(380, 5).Value = 154.12
$ws.Cells.Item(380, 6).Value = -56
$ws.Cells.Item(380, 7).Value = -7224.56

$ws.Cells.Item(382, 2).Value = 60325
$ws.Cells.Item(382, 5).Value = 151.57
$ws.Cells.Item(382, 6).Value = -102
$ws.Cells.Item(382, 7).Value = -12939.72

$ws.Cells.Item(383, 2).Value = 63560
$ws.Cells.Item(383, 5).Value = 134.87
$ws.Cells.Item(383, 6).Value = 104
$ws.Cells.Item(383, 7).Value = 13193.44

$ws.Cells.Item(400, 2).Value = 57835
$ws.Cells.Item(400, 6).Value = 1
$ws.Cells.Item(400, 7).Value = 59.13

$ws.Cells.Item(401, 2).Value = 62933
$ws.Cells.Item(401, 6).Value = 146
$ws.Cells.Item(401, 7).Value = 8632.98

$ws.Cells.Item(421, 2).Value = 57857
$ws.Cells.Item(421, 6).Value = 3
$ws.Cells.Item(421, 7).Value = 453.51

$ws.Cells.Item(422, 2).Value = 63008
$ws.Cells.Item(422, 6).Value = 504
$ws.Cells.Item(422, 7).Value = 76189.67999999999

$ws.Cells.Item(457, 2).Value = 31930
$ws.Cells.Item(457, 5).Value = 26.8
$ws.Cells.Item(457, 6).Value = -62
$ws.Cells.Item(457, 7).Value = -1390.04

$ws.Cells.Item(458, 2).Value = 63681
$ws.Cells.Item(458, 5).Value = 23.84
$ws.Cells.Item(458, 6).Value = 65
$ws.Cells.Item(458, 7).Value = 1457.3

$ws.Cells.Item(583, 2).Value = 53263
$ws.Cells.Item(583, 5).Value = 15.29
$ws.Cells.Item(583, 6).Value = -309
$ws.Cells.Item(583, 7).Value = -3958.29

$ws.Cells.Item(584, 2).Value = 65066
$ws.Cells.Item(584, 5).Value = 13.61
$ws.Cells.Item(584, 6).Value = 313
$ws.Cells.Item(584, 7).Value = 4009.53

$ws.Cells.Item(586, 2).Value = 64915
$ws.Cells.Item(586, 5).Value = 20.98
$ws.Cells.Item(586, 6).Value = 40
$ws.Cells.Item(586, 7).Value = 789.2

$ws.Cells.Item(587, 2).Value = 45695
$ws.Cells.Item(587, 5).Value = 23.58
$ws.Cells.Item(587, 6).Value = -36
$ws.Cells.Item(587, 7).Value = -710.28

$ws.Cells.Item(599, 2).Value = 64925
$ws.Cells.Item(599, 5).Value = 13.97
$ws.Cells.Item(599, 6).Value = 302
$ws.Cells.Item(599, 7).Value = 3971.3

$ws.Cells.Item(600, 2).Value = 45709
$ws.Cells.Item(600, 5).Value = 15.69
$ws.Cells.Item(600, 6).Value = -300
$ws.Cells.Item(600, 7).Value = -3945

$ws.Cells.Item(687, 2).Value = 64810
$ws.Cells.Item(687, 5).Value = 291.22
$ws.Cells.Item(687, 6).Value = 7
$ws.Cells.Item(687, 7).Value = 1917.44

$ws.Cells.Item(688, 2).Value = 53319
$ws.Cells.Item(688, 5).Value = 310.64
$ws.Cells.Item(688, 6).Value = -6
$ws.Cells.Item(688, 7).Value = -1643.52

$ws.Cells.Item(709, 2).Value = 64833
$ws.Cells.Item(709, 5).Value = 34.9
$ws.Cells.Item(709, 6).Value = 99
$ws.Cells.Item(709, 7).Value = 3250.17

$ws.Cells.Item(710, 2).Value = 60025
$ws.Cells.Item(710, 5).Value = 37.22
$ws.Cells.Item(710, 6).Value = -98
$ws.Cells.Item(710, 7).Value = -3217.34

$ws.Cells.Item(720, 2).Value = 64830
$ws.Cells.Item(720, 5).Value = 34.9
$ws.Cells.Item(720, 6).Value = 117
$ws.Cells.Item(720, 7).Value = 3841.11

$ws.Cells.Item(721, 2).Value = 60022
$ws.Cells.Item(721, 5).Value = 37.22
$ws.Cells.Item(721, 6).Value = -113
$ws.Cells.Item(721, 7).Value = -3709.79

$ws.Cells.Item(872, 2).Value = 65362
$ws.Cells.Item(872, 6).Value = 2
$ws.Cells.Item(872, 7).Value = 81.73999999999999

$ws.Cells.Item(873, 2).Value = 65079
$ws.Cells.Item(873, 6).Value = 21
$ws.Cells.Item(873, 7).Value = 858.27
